# Update the division problem answers in the table to match the new
# output generated at 9a8706d.

$d = $word.ActiveDocument

$replacements = @(
    @{old = "445÷3=148, 1"; new = "190÷3=63, 1"},
    @{old = "874÷6=145, 4"; new = "556÷8=69, 4"},
    @{old = "153÷6=25, 3"; new = "682÷4=170, 2"},
    @{old = "236÷8=29, 4"; new = "623÷8=77, 7"},
    @{old = "984÷2=492, 0"; new = "546÷3=182, 0"},
    @{old = "291÷4=72, 3"; new = "563÷9=62, 5"},
    @{old = "523÷8=65, 3"; new = "480÷6=80, 0"},
    @{old = "552÷5=110, 2"; new = "107÷7=15, 2"},
    @{old = "236÷9=26, 2"; new = "115÷9=12, 7"},
    @{old = "907÷6=151, 1"; new = "307÷9=34, 1"},
    @{old = "657÷2=328, 1"; new = "990÷3=330, 0"},
    @{old = "941÷8=117, 5"; new = "744÷6=124, 0"},
    @{old = "501÷9=55, 6"; new = "518÷4=129, 2"},
    @{old = "531÷9=59, 0"; new = "126÷8=15, 6"},
    @{old = "800÷7=114, 2"; new = "144÷3=48, 0"},
    @{old = "346÷8=43, 2"; new = "430÷2=215, 0"},
    @{old = "372÷7=53, 1"; new = "979÷9=108, 7"},
    @{old = "587÷8=73, 3"; new = "682÷2=341, 0"},
    @{old = "236÷3=78, 2"; new = "327÷6=54, 3"},
    @{old = "973÷8=121, 5"; new = "769÷5=153, 4"},
    @{old = "270÷8=33, 6"; new = "522÷8=65, 2"},
    @{old = "566÷3=188, 2"; new = "715÷7=102, 1"},
    @{old = "169÷3=56, 1"; new = "777÷2=388, 1"},
    @{old = "797÷8=99, 5"; new = "119÷7=17, 0"},
    @{old = "287÷5=57, 2"; new = "885÷9=98, 3"}
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $true, $false, $false, $false, `
                         $true, 1, $false, $pair.new, 2)
}

$d.Save()
